# Apply updated dSF (column F) values to specific rows, per the
# "repull data, push all data, mean calculation" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 0
    6  = -2
    9  = 5
    10 = 1
    11 = -5
    13 = 3
    15 = 2
    18 = -2
    20 = 0
    21 = 7
    28 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
